$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(97, 27, 27.00000000000221, 27, 27, 19, 20, 30.00000000001455)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
